# Resync the calendar availability sheet from the (authoritative) Excel input.
# 2026-02-25 and 2026-02-26 move from "Closed"/unavailable to "Limited"
# availability (a price reappears and one of the two bookable slots is now
# booked), and the stale trailing blank rows (7-9) are trimmed away.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: 2026-02-25 --------------------------------------------------
$ws.Range("B4").Value = "Limited"
$ws.Range("C4").Value = 18500
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

# --- Row 5: 2026-02-26 --------------------------------------------------
$ws.Range("B5").Value = "Limited"
$ws.Range("C5").Value = 18500
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1

# --- Drop the now-unused trailing blank rows -----------------------------
$ws.Rows("7:9").Delete()

# Re-measure rows 1-6 so the stale explicit row height is dropped in favor
# of the sheet's (now taller) default line height.
$ws.Rows("1:6").AutoFit()

# Leave the selection where the editor last clicked.
$null = $ws.Range("F11").Select()
